$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix test case 5: Zip Code values of "11616" should be "100000"
# Keep them as text (matching original text-typed cells) by using a
# leading apostrophe so Excel doesn't reinterpret them as numbers.
$ws.Range("E2").Value = "'100000"
$ws.Range("E3").Value = "'100000"
$ws.Range("E5").Value = "'100000"
$ws.Range("E6").Value = "'100000"

# Row 2 no longer needs its custom (header-like) row height; restore it
# to the sheet's default row height.
$ws.Rows(2).AutoFit()

# Remove the extra (duplicate) test-case row 7 entirely
$ws.Rows(7).Delete()

# Update the active selection to match the post-edit state
$ws.Range("D12").Select()
